# Apply roster reshuffle to Sheet1 (A2:C16)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Donovan Clingan", "C", "Portland Trail Blazers"),
    @("Buddy Hield", "SG,SF", "Golden State Warriors"),
    @("Onyeka Okongwu", "PF,C", "Atlanta Hawks"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
